# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit:
#  - Remove the second worker's (22803405 / HEYDI PATRICIA CASTRO GARCIA) 8 data rows.
#  - Re-sort the remaining worker's period rows into ascending period order.
#  - Update the "VALOR MORA" total and "Cant. Trabajadores" count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 8 rows that hold the second worker's period/mora data
# (old rows 25-32). This shifts the signature block (old rows 37-38) up to
# rows 29-30 and removes the now-unused shared strings automatically.
$ws.Range("B25:J32").EntireRow.Delete()

# Re-sort the remaining worker's rows (now B16:J24) into ascending period
# order (2011, 2012, 2101..2107), keeping each period's original Valor Mora.
$periods = @("2011", "2012", "2101", "2102", "2103", "2104", "2105", "2106", "2107")
$valores = @(40000, 35112, 35112, 35112, 35112, 35112, 35112, 35112, 29260)

for ($i = 0; $i -lt 9; $i++) {
    $row = 16 + $i
    $ws.Range("E" + $row).Value = $periods[$i]
    $ws.Range("F" + $row).Value = $valores[$i]
}

# Update totals: Valor Mora total and Cant. Trabajadores
$ws.Range("E11").Value = 315044
$ws.Range("C13").Value = 1
